# ============================================================================
# Edit script: dmas_solr_to_schema_dot_org_mapping.xlsx
#
# Adds a "Reference" (column D) to the DMAS Solr -> schema.org mapping
# table, pointing each mapping row at the relevant science-on-schema.org
# guide section, and appends new mapping rows for:
#   temporalCoverage, spatialCoverage, additionalProperty (ellipsoid/CRS),
#   creator, provider, publisher, and DatasetFunding.
# Also adds a live hyperlink on the "variableMeasured" reference cell.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: Add Reference column (D) values for existing rows 1-24 (except row 21) ----
$ws.Range("D1").Value = 'Reference'
$ws.Range("D2").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D3").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D4").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D5").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D6").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D7").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D8").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D9").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D10").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-variables'
$ws.Range("D11").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-variables'
$ws.Range("D12").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-variables'
$ws.Range("D13").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D14").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D15").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D16").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset'
$ws.Range("D17").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset-identifier'
$ws.Range("D18").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset-identifier'
$ws.Range("D19").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset-identifier'
$ws.Range("D20").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-dataset-identifier'
$ws.Range("D22").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-distributions'
$ws.Range("D23").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-distributions'
$ws.Range("D24").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-distributions'

# ---- Step 2: Add new rows 25-31 (temporal/spatial coverage, creator/provider/publisher, etc.) ----
$ws.Range("A25").Value = 'Dataset-DatasetCoverage-StartTimeLong + "/" + Dataset-DatasetCoverage-StopTimeLong'
$ws.Range("B25").Value = 'temporalCoverage'
$ws.Range("C25").Value = 'Temporal coverage at PO.DAAC is always defined using a date range e.g. 2012-09-20 - 2016-01-22 for example. The datetime needs to be mapped from the Long datetime included within the DMAS Solr response to the ISO 8601 equivalent. Additionally, if the end time is present or ongoing, then the datetime range can be expressed as follows "2012-09-20/.." note the two dots in the end date entry. This is documented further at https://github.com/schemaorg/schemaorg/issues/242. Also note, that the DMAS Solr response includes three varieties of start time e.g. DatasetCoverage-StartTimeLong-Long and DatasetCoverage-StartTimeLong in addition to Dataset-DatasetCoverage-StartTimeLong. The same is true for end times.'
$ws.Range("D25").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-temporal-coverage'

$ws.Range("A26").Value = 'DatasetCoverage-NorthLat, DatasetCoverage-SouthLat, DatasetCoverage-WestLon, DatasetCoverage-EastLon'
$ws.Range("B26").Value = 'spatialCoverage'
$ws.Range("C26").Value = 'This is always represented as a schema.org/Place of type GeoShape of type ''box''. See the accompanying Dataset.jsonld document for an example. In the example provided the ''lower-left'' corner is 39.3280/120.1633 and ''upper-right'' corner is 40.445/123.7878'
$ws.Range("D26").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-spatial-coverage'

$ws.Range("A27").Value = 'Dataset-EllipsoidType'
$ws.Range("B27").Value = 'additionalProperty'
$ws.Range("C27").Value = 'See the example for how one would map a value to CRS84. Lot''s of the PO.DAAC Datasets have WGS84 ellipsoid representations however, so we need to map those differently. See both the science-on-schema.org link and the Dataset.json example.'
$ws.Range("D27").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#spatial_multiple-geometries'

$ws.Range("B28").Value = 'creator'
$ws.Range("C28").Value = 'Descrbing a datasets people is not particularly of interest to PO.DAAC. The reason here is that we do not necessarily wish for people included in the DMAS Solr response to be stated as points of contact for the PO.DAAC dataset.'
$ws.Range("D28").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-people'

$ws.Range("A29").Value = 'N/A'
$ws.Range("B29").Value = 'provider'
$ws.Range("C29").Value = 'Same as for DataRepository, this value is is simply ''https://podaac.jpl.nasa.gov'''
$ws.Range("D29").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-publisherprovider'

$ws.Range("A30").Value = 'N/A'
$ws.Range("B30").Value = 'publisher'
$ws.Range("C30").Value = 'Same as for DataRepository, this value is is simply ''https://podaac.jpl.nasa.gov'''
$ws.Range("D30").Value = 'https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md#describing-a-datasets-publisherprovider'

$ws.Range("B31").Value = 'DatasetFunding'
$ws.Range("C31").Value = 'Same as with ''creator'' we most likely do NOT wish to go into details of who funded the dataset. Right now, I don''t think we host that kind of information at PO.DAAC right now anyway. This is another open question however.'

# ---- Step 3: Turn the D10 reference into a real hyperlink (cell text is
#      already the full URL with fragment from Step 1; Hyperlinks.Add below
#      only attaches the link + applies the built-in "Hyperlink" style,
#      since we don't pass a TextToDisplay argument). ----
$ws.Hyperlinks.Add($ws.Range("D10"), "https://github.com/ESIPFed/science-on-schema.org/blob/master/guides/Dataset.md", "describing-a-datasets-variables")

# ---- Step 4: Match row heights to the new wrapped content (auto-fit
#      equivalent, since wrap height can't be measured headlessly). ----
$ws.Rows.Item(25).RowHeight = 119
$ws.Rows.Item(26).RowHeight = 51
$ws.Rows.Item(27).RowHeight = 51
$ws.Rows.Item(28).RowHeight = 51
$ws.Rows.Item(29).RowHeight = 17
$ws.Rows.Item(30).RowHeight = 17
$ws.Rows.Item(31).RowHeight = 51

# ---- Step 5: Restore the view to roughly where the author left off. ----
$ws.Range("B33").Select()
